# Generate Report for Handback
#
# The handback has now completed for both files/languages, so:
#   - the overall Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown
#     (Overview sheet + both language sheets)
#   - the "Latest Target File" / "Latest Handback File" columns (F/G) on
#     each language sheet get populated (with hyperlinks) for both rows
#   - the "Latest Handback DateTime" column (H) on each language sheet
#     gets the real handback timestamp instead of the zero-date
#     placeholder

$wb = $excel.ActiveWorkbook

$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$mdFile2 = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/5c6331c8b284af3781cf5ffddca07ac2729da0b7/e2e/0055612a-0e2e-4e43-b43c-dc47a2018e97.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/5c6331c8b284af3781cf5ffddca07ac2729da0b7/e2e/bed71e79-5634-4f78-be6c-f2c1dca6cf5b.md"

# Cornflower-blue underlined hyperlink look used by the existing A/D columns.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status cells (B/C columns).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf1 = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
$zhXlf2 = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
$zhHandbackUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a95fab13b9e6cd842951b26294649a9dfee2101/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.zh-cn.xlf"
$zhHandbackUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a95fab13b9e6cd842951b26294649a9dfee2101/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.zh-cn.xlf"
$zhHandbackDateTime = "2016-03-30 11:00:46"

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("F2").Value = $mdFile1
$wsZh.Range("G2").Value = $zhXlf1
$wsZh.Range("H2").Value = $zhHandbackDateTime

$wsZh.Range("F3").Value = $mdFile2
$wsZh.Range("G3").Value = $zhXlf2
$wsZh.Range("H3").Value = $zhHandbackDateTime

# Only touch the brand-new cells' hyperlinks -- leave the pre-existing
# A2/D2/A3/D3 hyperlinks (and their styling) completely alone.
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhHandbackUrl1, "", "", $zhXlf1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhHandbackUrl2, "", "", $zhXlf2) | Out-Null

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $wsZh.Range($addr).Font.Underline = 2
    $wsZh.Range($addr).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf1 = "0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
$deXlf2 = "bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
$deHandbackUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/857c1a242ce9404ed0983a34e623258443a7c4cf/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/0055612a-0e2e-4e43-b43c-dc47a2018e97.3cb6669bbc1d8c1478dd6c22dae19bb2cb345e54.de-de.xlf"
$deHandbackUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/857c1a242ce9404ed0983a34e623258443a7c4cf/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/bed71e79-5634-4f78-be6c-f2c1dca6cf5b.a3205878e3a5027496ef44cb964c4bde2ca4dc1b.de-de.xlf"
$deHandbackDateTime = "2016-03-30 11:01:05"

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("F2").Value = $mdFile1
$wsDe.Range("G2").Value = $deXlf1
$wsDe.Range("H2").Value = $deHandbackDateTime

$wsDe.Range("F3").Value = $mdFile2
$wsDe.Range("G3").Value = $deXlf2
$wsDe.Range("H3").Value = $deHandbackDateTime

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl1, "", "", $mdFile1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deHandbackUrl1, "", "", $deXlf1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl2, "", "", $mdFile2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deHandbackUrl2, "", "", $deXlf2) | Out-Null

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $wsDe.Range($addr).Font.Underline = 2
    $wsDe.Range($addr).Font.Color = $hyperlinkColor
}
